$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attack")

# Add the description text for the default passive in C2
$ws.Range("C2").Value = "This is the default passive. Straight up you don't want this it does nuttin."

# Update the active selection to C2 (matches the saved selection state in the diff)
$ws.Range("C2").Select()
